# "Update : multiple file export by sheet name"
#
# - Sheet1 ("Sheet1") is renamed to "app"
# - Sheet2 ("Sheet2") is renamed to "other" and filled in with a copy of
#   the "app" key/translation table, but describing an "OTHER*" key set
# - "other" becomes the active tab/sheet
# - Sheet1's own selection is cleared back to the whole table, "other"'s
#   selection sits just below its data (B10)

$wb  = $excel.ActiveWorkbook
$app = $wb.Worksheets.Item(1)
$other = $wb.Worksheets.Item(2)

$app.Name = "app"
$other.Name = "other"

# ---- Clone app's A1:D7 layout (values + formatting + merges) into other ----
$app.Range("A1:D7").Copy()
$other.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$other.Range("A1").PasteSpecial(-4163)   # xlPasteValues

$other.Range("A1:A2").MergeCells = $true
$other.Range("B1:D1").MergeCells = $true

# Re-apply the source formats on top of the merge so the merged cells keep
# the exact same style ids as "app" rather than the ones the merge step
# synthesizes for the (now split) borders.
$app.Range("A1:D7").Copy()
$other.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

$other.Range("A1:D7").RowHeight = 15.75

# ---- Column widths for "other" (narrower key column, wider lang columns) ----
$other.Columns.Item(1).ColumnWidth = 9.8
$other.Range("B1:D1").ColumnWidth = 26

# ---- Overwrite the key / translation cells with the "OTHER" data set ----
$other.Range("A3").Value = "OTHER1"

$other.Range("A4").Value = "OTHER2"
$other.Range("B4").Value = "EN OTHER 2"
$other.Range("C4").Value = "ID OTHER 2"
$other.Range("D4").Value = "AR OTHER 2"

$other.Range("A5").Value = "OTHER3"
$other.Range("B5").Value = "EN OTHER 3"
$other.Range("C5").Value = "ID OTHER 3"
$other.Range("D5").Value = "AR OTHER 3"

$other.Range("A6").Value = "OTHER4"
$other.Range("B6").Value = "EN OTHER 4"
$other.Range("C6").Value = "ID OTHER 4"
$other.Range("D6").Value = "AR OTHER 4"

$other.Range("A7").Value = "OTHER5"
$other.Range("B7").Value = "EN OTHER 5"
$other.Range("C7").Value = "ID OTHER 5"
$other.Range("D7").Value = "AR OTHER 5"

# ---- Print setup on "other" (matches "app"'s portrait page setup) ----
$other.PageSetup.Orientation = 1   # xlPortrait

# ---- Selections / active tab ----
$app.Range("A1:D7").Select()
$other.Range("B10").Select()
$other.Activate()

Write-Host "done"
